# plasmids.xlsx: add plasmid "features" (Resistance / Origin) columns.
#
# The sheet used to be a single "Name" column (A1 header, A2 = "p2").
# We now want a 2-column table describing plasmid features instead:
#   A1="Resistance"  B1="Origin"
#   A3="AmpR,TetR"   B3="p15A"
# (row 2 is intentionally left blank, matching the target layout).
#
# Deleting column A first (instead of just overwriting the old cells)
# is what naturally produces the shifted column widths seen in the
# target file - every remaining column's width slides one slot to the
# left, exactly like Excel does on a real "Delete Entire Column".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old single "Name" column entirely (shifts B..G left to A..F,
# carrying their widths with them) and drops the old "Name"/"p2" values.
$ws.Range("A1").EntireColumn.Delete()

# New header row.
$ws.Range("A1").Value = "Resistance"
$ws.Range("B1").Value = "Origin"

# New data row (row 2 stays empty on purpose).
$ws.Range("A3").Value = "AmpR,TetR"
$ws.Range("B3").Value = "p15A"

# Park the selection back on A1 (was C6 previously).
[void]$ws.Range("A1").Select()
